# update coupon code for MWL
# Sheet1!A2 holds the coupon/id code for this row; replace the old code
# ("CA-MWQYTQLX") with the newly issued one ("CA-XGYTNHX3").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").Value = "CA-XGYTNHX3"
